# Consolidate the duplicate "UEMC" team rows (row 29 = "UEMC BALONCESTO
# VALLADOLID", row 30 = "UEMC CBC VALLADOLID") into a single row 29, then
# remove the now-redundant row 30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ratio / rate columns (C..J): consolidated value is the average of the two
# rows' values.
$ratioCols = "C","D","E","F","G","H","I","J"
foreach ($col in $ratioCols) {
    $a = $ws.Range("$col`29").Value2
    $b = $ws.Range("$col`30").Value2
    $ws.Range("$col`29").Value2 = ($a + $b) / 2
}

# Count / total columns (K..AC): consolidated value is the sum of the two
# rows' values.
$sumCols = "K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC"
foreach ($col in $sumCols) {
    $a = $ws.Range("$col`29").Value2
    $b = $ws.Range("$col`30").Value2
    $ws.Range("$col`29").Value2 = $a + $b
}

# Remove row 30 (now merged into row 29) entirely, shrinking the sheet's
# used range back down to A1:AC29.
$ws.Rows.Item(30).Delete()
